$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Use git to clone the source." code sample paragraph:
#    The literal text runs were split as "g" | <bookmark _GoBack> | "it clone "
#    | "ghostscript.com:/home/fred/repos/qt-gsview.git". Word had stashed the
#    "_GoBack" bookmark in the middle of the word "git" (a relic of the last
#    edit position). We want the paragraph to read as a clean "git clone "
#    run (merged) followed by the URL run, with no bookmark splitting "git".
$gitPara = $d.Paragraphs.Item(19)
$gitRange = $gitPara.Range

# Temporarily mark the boundary between "it clone " and "ghostscript..." so
# that when we touch/merge the "g"/"it clone " runs below, this boundary
# (which must remain two distinct runs) is not coalesced away too.
$mergeGuard = $d.Bookmarks.Add("ZZGuard", $d.Range($gitRange.Start + 10, $gitRange.Start + 10))

# Drop the stray _GoBack bookmark that was sitting between "g" and "it clone ".
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Touch the text spanning the old "g" | "it clone " run boundary so the two
# (now bookmark-free, identically formatted) runs coalesce into one run.
$gitRange.Find.Execute("git clone", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "git clone", 2) | Out-Null

# Clean up the temporary guard bookmark.
$mergeGuard2 = $d.Bookmarks.Item("ZZGuard")
$mergeGuard2.Delete()

# ---------------------------------------------------------------------
# 2) Both "./configure --with-libiconv=no" build-step lines gain the
#    "--disable-cups" flag.
$d.Content.Find.Execute("./configure --with-libiconv=no", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "./configure --with-libiconv=no --disable-cups", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) The cursor's last-known position (_GoBack) now belongs on the empty
#    "make" paragraph right after the second "./configure" line.
$makePara = $d.Paragraphs.Item(44)
$makeStart = $makePara.Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($makeStart, $makeStart)) | Out-Null
